$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "WIP_testing predefined models\BESOS_Output\X4E19VWN0PTA2159A5GZ"
$ws.Range("M3").Value = "WIP_testing predefined models\BESOS_Output\A6SWO6FY6Z7J0MVGCRJK"
$ws.Range("M4").Value = "WIP_testing predefined models\BESOS_Output\YQJC31UPERDY6A8NBE84"
$ws.Range("M5").Value = "WIP_testing predefined models\BESOS_Output\SJ47E8CQWB3RUHK3GJ6I"
